# Generate Report for Handback
# Updates the "Status" and "Error Detail" columns for the row corresponding to
# 5e9e483d-6369-45d7-a849-d54b42e65dda on both the "zh-cn" and "de-de" sheets,
# and widens the "Error Detail" column so the new message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Widen the "Error Detail" column (column P) on both locale sheets.
# (39.15 is the COM ColumnWidth input that round-trips to a saved raw
# column width of exactly 40 in the underlying OOXML.)
$zhcn.Columns.Item(16).ColumnWidth = 39.15
$dede.Columns.Item(16).ColumnWidth = 39.15

# Row 3 on each sheet corresponds to file 5e9e483d-6369-45d7-a849-d54b42e65dda.
# Update the Status column (C) on both locale sheets, and the matching
# per-locale status mirror columns (E: zh-cn, F: de-de) on the Overview sheet.
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# Populate the Error Detail column (P) with the handback failure reason.
$zhcn.Range("P3").Value = "Handback file name: qhxxhdjl.wxe is different with handoff file name: 5e9e483d-6369-45d7-a849-d54b42e65dda.5cb8c4b3c277fb3e49b163f5e08771ff90b72840.zh-cn."
$dede.Range("P3").Value = "Handback file name: qhxxhdjl.wxe is different with handoff file name: 5e9e483d-6369-45d7-a849-d54b42e65dda.5cb8c4b3c277fb3e49b163f5e08771ff90b72840.de-de."
